# Generate Report for Handoff
# b.md has been handed off again (new zh-cn / de-de target files produced).
# Update the Overview sheet's status for b.md, and the per-locale detail
# sheets (zh-cn, de-de) with the new "Latest Handoff File" / "Latest Handoff
# Datetime" for b.md's row, including the hyperlink label on that cell.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet: row for b.md (row 3) -> Status columns (B, C) change
# from "Handed back: in sync with en-US" to "Ready for handoff"
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = "Ready for handoff"
$overview.Range("C3").Value = "Ready for handoff"

# ---------------------------------------------------------------------
# zh-cn sheet: row for b.md (row 3)
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("B3").Value = "Ready for handoff"
$zhcn.Range("C3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("D3").Value = "2016-03-11 04:47:06"

# Refresh the hyperlinks on this sheet so the C3 link's display text
# matches the new handoff file name (target URL is unchanged).
$zhcn.Hyperlinks.Delete()
$zhcn.Hyperlinks.Add($zhcn.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/1f8a40961d69b67c5130d189ba60d41011c0d677/e2e/a.md", "", "", "a.md")
$zhcn.Hyperlinks.Add($zhcn.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ac5c37d8c2e13e594ed518adfbd731f751b8040d/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf", "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf")
$zhcn.Hyperlinks.Add($zhcn.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/903dbb281af0c114a42eab7c17f12d5a1f70c7a0/e2e/a.md", "", "", "a.md")
$zhcn.Hyperlinks.Add($zhcn.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/355adc080f5f816ec0f55b9ca21a3e09db5bf9eb/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf", "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf")
$zhcn.Hyperlinks.Add($zhcn.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/1f8a40961d69b67c5130d189ba60d41011c0d677/e2e/b.md", "", "", "b.md")
$zhcn.Hyperlinks.Add($zhcn.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ac5c37d8c2e13e594ed518adfbd731f751b8040d/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf", "", "", "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf")
$zhcn.Hyperlinks.Add($zhcn.Range("E3"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/903dbb281af0c114a42eab7c17f12d5a1f70c7a0/e2e/a.md", "", "", "a.md")
$zhcn.Hyperlinks.Add($zhcn.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/355adc080f5f816ec0f55b9ca21a3e09db5bf9eb/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf", "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf")
$zhcn.Hyperlinks.Add($zhcn.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/1f8a40961d69b67c5130d189ba60d41011c0d677/.localization-config", "", "", ".localization-config")

# ---------------------------------------------------------------------
# de-de sheet: row for b.md (row 3)
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("B3").Value = "Ready for handoff"
$dede.Range("C3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("D3").Value = "2016-03-11 04:47:14"

# Refresh the hyperlinks on this sheet so the C3 link's display text
# matches the new handoff file name (target URL is unchanged).
$dede.Hyperlinks.Delete()
$dede.Hyperlinks.Add($dede.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/1f8a40961d69b67c5130d189ba60d41011c0d677/e2e/a.md", "", "", "a.md")
$dede.Hyperlinks.Add($dede.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/97f0678b44e0248e7c58633ee152934ebe2e0301/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf", "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf")
$dede.Hyperlinks.Add($dede.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/0a7de43e2c756ac2aadfd0335219227e71391154/e2e/a.md", "", "", "a.md")
$dede.Hyperlinks.Add($dede.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/6ff5574df427f1e5c8f4b9b96073900627fdb548/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf", "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf")
$dede.Hyperlinks.Add($dede.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/1f8a40961d69b67c5130d189ba60d41011c0d677/e2e/b.md", "", "", "b.md")
$dede.Hyperlinks.Add($dede.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/97f0678b44e0248e7c58633ee152934ebe2e0301/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf", "", "", "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf")
$dede.Hyperlinks.Add($dede.Range("E3"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/0a7de43e2c756ac2aadfd0335219227e71391154/e2e/a.md", "", "", "a.md")
$dede.Hyperlinks.Add($dede.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/6ff5574df427f1e5c8f4b9b96073900627fdb548/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf", "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf")
$dede.Hyperlinks.Add($dede.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/1f8a40961d69b67c5130d189ba60d41011c0d677/.localization-config", "", "", ".localization-config")

Write-Output "Report generated for handoff of b.md"
